$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adjust the "weight" (module) on row 16: change B16 from WEB1201 to NET3204
$ws.Range("B16").Value = "NET3204 "

# Improve timetable exporter: clear out semester 5 rows (18-21), leaving the
# formatted-but-empty cells behind
$ws.Range("A18:B21").ClearContents()

# Update the active selection to reflect the exported/reviewed range
$ws.Range("A2:B17").Select()
